$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = "6 mm"
$ws.Range("B7").Value = "13 mm + 1.5*tol"
$ws.Range("B22").Value = "2.25 mm"

$ws.Range("G25").Select()
